$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the "Load in necessary packages (...)" list item:
#    rvest -> tidyverse, chromote -> lubridate, and drop the now
#    duplicated ", tidyverse, lubridate" tail so the parenthetical
#    ends up as "(tidyverse, lubridate)".
# ------------------------------------------------------------------
$pPackages = $d.Paragraphs(4)
$pPackages.Range.Find.Execute("rvest", $false, $false, $false, $false, $false, `
    $true, 1, $false, "tidyverse", 2) | Out-Null
$pPackages.Range.Find.Execute("chromote", $false, $false, $false, $false, $false, `
    $true, 1, $false, "lubridate", 2) | Out-Null
$pPackages.Range.Find.Execute(", tidyverse, lubridate)", $false, $false, $false, $false, $false, `
    $true, 1, $false, ")", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the blank ListParagraph that followed the packages line.
# ------------------------------------------------------------------
$pBlank = $d.Paragraphs(5)
$pBlank.Range.Delete()

# ------------------------------------------------------------------
# 3) Insert a new "Read in Tremblant_Data.csv" bullet right before
#    the "Write a function that converts time to seconds" bullet.
# ------------------------------------------------------------------
$pSeconds = $d.Paragraphs(5)
$pSeconds.Range.InsertParagraphBefore()
$pReadCsv = $d.Paragraphs(5)
$pReadCsv.Range.InsertBefore("Read in Tremblant_Data.csv")

# ------------------------------------------------------------------
# 4) Delete the old scraping-function block in its entirety: the two
#    blank paragraphs, the "Write a function to scrape..." intro,
#    the seven numId=2 sub-bullets, and the trailing ind-1080 blank
#    paragraph. That is paragraphs 7 through 17 inclusive.
# ------------------------------------------------------------------
$pBlockStart = $d.Paragraphs(7)
$pBlockEnd = $d.Paragraphs(17)
$blockRange = $d.Range($pBlockStart.Range.Start, $pBlockEnd.Range.End)
$blockRange.Delete()

# ------------------------------------------------------------------
# 5) Append the new cleaning-module bullets (numId=1) after the
#    "Test out the function with an id code" bullet, then replace
#    that bullet's text with "Fix Total using the function above".
# ------------------------------------------------------------------
$pTestCode = $d.Paragraphs(7)
$newBulletTexts = @(
    "Separate Run 1 into two columns for Run 1 time and Run 1 rank",
    "Follow similar steps as above for Run 2",
    "Use a function to make sure there are only numeric values in the Pr, Run 1 Rank, and Run 2 Rank columns",
    "Use your function to fix time on Run 1 time and Run 2 time",
    "Rename variables logically",
    "Order the variables logically"
)
$anchor = $pTestCode
foreach ($bulletText in $newBulletTexts) {
    $anchor.Range.InsertParagraphAfter()
    $insertedIndex = $anchor.Index + 1
    $insertedPara = $d.Paragraphs($insertedIndex)
    $insertedPara.Range.InsertBefore($bulletText)
    $anchor = $insertedPara
}

$pTestCode = $d.Paragraphs(7)
$pTestCode.Range.Delete()

$pFixTotalSlot = $d.Paragraphs(7)
$pFixTotalSlot.Range.InsertParagraphBefore()
$pFixTotal = $d.Paragraphs(7)
$pFixTotal.Range.InsertBefore("Fix Total using the function above")
